$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update Valor Mora (E11)
$ws.Range("E11").Value = 106000

# Update Cant. Periodos (F13)
$ws.Range("F13").Value = 1

# Delete rows 17 and 18 (the 2501 and 2412 period rows), shifting rows below up
$ws.Range("A17:A18").EntireRow.Delete()
